$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.0275626423690205
$ws.Range("J2").Value = 0.02947285804177352
$ws.Range("M2").Value = 8.424749
$ws.Range("N2").Value = 16.849498
$ws.Range("O2").Value = 0.1980372819757593
$ws.Range("P2").Value = 0.1585662899160533
$ws.Range("Q2").Value = 0.012232735548
$ws.Range("R2").Value = 0.07339641328800001
$ws.Range("S2").Value = 0.005458430778830725
$ws.Range("T2").Value = 0.004673401752906542
$ws.Range("I3").Value = 0.0275626423690205
$ws.Range("J3").Value = 0.02947285804177352
$ws.Range("O3").Value = 0.2481739331426511
$ws.Range("P3").Value = 0.2980652388254983
$ws.Range("S3").Value = 0.006840329364524098
$ws.Range("T3").Value = 0.008784834471091232
$ws.Range("I4").Value = 0.0275626423690205
$ws.Range("J4").Value = 0.02947285804177352
$ws.Range("M4").Value = 5.743874333333333
$ws.Range("N4").Value = 17.231623
$ws.Range("O4").Value = 0.1350190089916815
$ws.Range("P4").Value = 0.1621623699615343
$ws.Range("Q4").Value = 0.008340105532
$ws.Range("R4").Value = 0.075060949788
$ws.Range("S4").Value = 0.003721480657857282
$ws.Range("T4").Value = 0.004779388509593857
$ws.Range("I5").Value = 0.0275626423690205
$ws.Range("J5").Value = 0.02947285804177352
$ws.Range("M5").Value = 12.9373935
$ws.Range("N5").Value = 25.874787
$ws.Range("O5").Value = 0.3041142524947457
$ws.Range("P5").Value = 0.2435009622813764
$ws.Range("Q5").Value = 0.018785095362
$ws.Range("R5").Value = 0.112710572172
$ws.Range("S5").Value = 0.008382192380834676
$ws.Range("T5").Value = 0.007176669294354253
$ws.Range("I6").Value = 0.0275626423690205
$ws.Range("J6").Value = 0.02947285804177352
$ws.Range("M6").Value = 0.882742
$ws.Range("N6").Value = 2.648226
$ws.Range("O6").Value = 0.02075027117909931
$ws.Range("P6").Value = 0.02492177343676531
$ws.Range("Q6").Value = 0.001281741384
$ws.Range("R6").Value = 0.011535672456
$ws.Range("S6").Value = 0.0005719323035697077
$ws.Range("T6").Value = 0.000734515890651026
$ws.Range("I7").Value = 0.0275626423690205
$ws.Range("J7").Value = 0.02947285804177352
$ws.Range("M7").Value = 3.994844666666667
$ws.Range("N7").Value = 11.984534
$ws.Range("O7").Value = 0.09390525221606305
$ws.Range("P7").Value = 0.1127833655787726
$ws.Range("Q7").Value = 0.005800514456000001
$ws.Range("R7").Value = 0.052204630104
$ws.Range("S7").Value = 0.002588276883404016
$ws.Range("T7").Value = 0.00332404812317661
$ws.Range("I8").Value = 0.1944381169324222
$ws.Range("J8").Value = 0.1386090380724913
$ws.Range("M8").Value = 8.424749
$ws.Range("N8").Value = 16.849498
$ws.Range("O8").Value = 0.1980372819757593
$ws.Range("P8").Value = 0.1585662899160533
$ws.Range("Q8").Value = 0.08629470400700001
$ws.Range("R8").Value = 0.345178816028
$ws.Range("S8").Value = 0.03850599618978175
$ws.Range("T8").Value = 0.02197872091598792
$ws.Range("I9").Value = 0.1944381169324222
$ws.Range("J9").Value = 0.1386090380724913
$ws.Range("O9").Value = 0.2481739331426511
$ws.Range("P9").Value = 0.2980652388254983
$ws.Range("S9").Value = 0.04825447223196992
$ws.Range("T9").Value = 0.04131453603644972
$ws.Range("I10").Value = 0.1944381169324222
$ws.Range("J10").Value = 0.1386090380724913
$ws.Range("M10").Value = 5.743874333333333
$ws.Range("N10").Value = 17.231623
$ws.Range("O10").Value = 0.1350190089916815
$ws.Range("P10").Value = 0.1621623699615343
$ws.Range("Q10").Value = 0.05883450479633334
$ws.Range("R10").Value = 0.353007028778
$ws.Range("S10").Value = 0.02625284185842433
$ws.Range("T10").Value = 0.02247717011192372
$ws.Range("I11").Value = 0.1944381169324222
$ws.Range("J11").Value = 0.1386090380724913
$ws.Range("M11").Value = 12.9373935
$ws.Range("N11").Value = 25.874787
$ws.Range("O11").Value = 0.3041142524947457
$ws.Range("P11").Value = 0.2435009622813764
$ws.Range("Q11").Value = 0.1325177216205
$ws.Range("R11").Value = 0.530070886482
$ws.Range("S11").Value = 0.05913140258738952
$ws.Range("T11").Value = 0.03375143415154756
$ws.Range("I12").Value = 0.1944381169324222
$ws.Range("J12").Value = 0.1386090380724913
$ws.Range("M12").Value = 0.882742
$ws.Range("N12").Value = 2.648226
$ws.Range("O12").Value = 0.02075027117909931
$ws.Range("P12").Value = 0.02492177343676531
$ws.Range("Q12").Value = 0.009041926306
$ws.Range("R12").Value = 0.054251557836
$ws.Range("S12").Value = 0.004034643653901182
$ws.Range("T12").Value = 0.003454383043130605
$ws.Range("I13").Value = 0.1944381169324222
$ws.Range("J13").Value = 0.1386090380724913
$ws.Range("M13").Value = 3.994844666666667
$ws.Range("N13").Value = 11.984534
$ws.Range("O13").Value = 0.09390525221606305
$ws.Range("P13").Value = 0.1127833655787726
$ws.Range("Q13").Value = 0.04091919392066667
$ws.Range("R13").Value = 0.245515163524
$ws.Range("S13").Value = 0.01825876041095546
$ws.Range("T13").Value = 0.0156327938134518
$ws.Range("G14").Value = 0.040985
$ws.Range("H14").Value = 0.122955
$ws.Range("I14").Value = 0.7779992406985573
$ws.Range("J14").Value = 0.8319181038857351
$ws.Range("M14").Value = 8.424749
$ws.Range("N14").Value = 16.849498
$ws.Range("O14").Value = 0.1980372819757593
$ws.Range("P14").Value = 0.1585662899160533
$ws.Range("Q14").Value = 0.345288337765
$ws.Range("R14").Value = 2.07173002659
$ws.Range("S14").Value = 0.1540728550071468
$ws.Range("T14").Value = 0.1319141672471588
$ws.Range("G15").Value = 0.040985
$ws.Range("H15").Value = 0.122955
$ws.Range("I15").Value = 0.7779992406985573
$ws.Range("J15").Value = 0.8319181038857351
$ws.Range("O15").Value = 0.2481739331426511
$ws.Range("P15").Value = 0.2980652388254983
$ws.Range("Q15").Value = 0.4327042059783333
$ws.Range("R15").Value = 3.894337853805
$ws.Range("S15").Value = 0.1930791315461571
$ws.Range("T15").Value = 0.2479658683179574
$ws.Range("G16").Value = 0.040985
$ws.Range("H16").Value = 0.122955
$ws.Range("I16").Value = 0.7779992406985573
$ws.Range("J16").Value = 0.8319181038857351
$ws.Range("M16").Value = 5.743874333333333
$ws.Range("N16").Value = 17.231623
$ws.Range("O16").Value = 0.1350190089916815
$ws.Range("P16").Value = 0.1621623699615343
$ws.Range("Q16").Value = 0.2354126895516667
$ws.Range("R16").Value = 2.118714205965
$ws.Range("S16").Value = 0.1050446864753999
$ws.Range("T16").Value = 0.1349058113400167
$ws.Range("G17").Value = 0.040985
$ws.Range("H17").Value = 0.122955
$ws.Range("I17").Value = 0.7779992406985573
$ws.Range("J17").Value = 0.8319181038857351
$ws.Range("M17").Value = 12.9373935
$ws.Range("N17").Value = 25.874787
$ws.Range("O17").Value = 0.3041142524947457
$ws.Range("P17").Value = 0.2435009622813764
$ws.Range("Q17").Value = 0.5302390725974999
$ws.Range("R17").Value = 3.181434435584999
$ws.Range("S17").Value = 0.2366006575265215
$ws.Range("T17").Value = 0.2025728588354745
$ws.Range("G18").Value = 0.040985
$ws.Range("H18").Value = 0.122955
$ws.Range("I18").Value = 0.7779992406985573
$ws.Range("J18").Value = 0.8319181038857351
$ws.Range("M18").Value = 0.882742
$ws.Range("N18").Value = 2.648226
$ws.Range("O18").Value = 0.02075027117909931
$ws.Range("P18").Value = 0.02492177343676531
$ws.Range("Q18").Value = 0.03617918087
$ws.Range("R18").Value = 0.32561262783
$ws.Range("S18").Value = 0.01614369522162842
$ws.Range("T18").Value = 0.02073287450298367
$ws.Range("G19").Value = 0.040985
$ws.Range("H19").Value = 0.122955
$ws.Range("I19").Value = 0.7779992406985573
$ws.Range("J19").Value = 0.8319181038857351
$ws.Range("M19").Value = 3.994844666666667
$ws.Range("N19").Value = 11.984534
$ws.Range("O19").Value = 0.09390525221606305
$ws.Range("P19").Value = 0.1127833655787726
$ws.Range("Q19").Value = 0.1637287086633333
$ws.Range("R19").Value = 1.47355837797
$ws.Range("S19").Value = 0.07305821492170357
$ws.Range("T19").Value = 0.09382652364214418
